$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E1 was a numeric date (03/03/2024); replace with the literal text "03_03_2024"
# matching the style of the other header cells (B1, C1, D1).
$ws.Range("E1").Value = "03_03_2024"

# Update selection to E2 (as reflected in the saved view state).
$ws.Range("E2").Select()
